$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Rename product name on both sheets
$newProductName = "4229-RBI-EI-DB-SAR-REC-CTRFD-RNI-FEE+INTEREST-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Change shortname from numeric 4229 to text "422i"
$wsInput.Range("B2").Value = "422i"

# Update selection on the input sheet
$wsInput.Range("B8").Select()

# Make the output sheet the active tab
$wsOutput.Activate()
